# Rotate the species-observation data in rows 2-4 cyclically:
#   old row 2 -> new row 3
#   old row 3 -> new row 4
#   old row 4 -> new row 2
# Only columns A, B, E, F, G, H, Q, R, AO differ between the rows; every
# other column already holds identical values across the three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AO")

# Capture the current ("before") values for the columns that change.
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value()
    $row3[$col] = $ws.Range($col + "3").Value()
    $row4[$col] = $ws.Range($col + "4").Value()
}

# Write back the rotated values.
foreach ($col in $cols) {
    $ws.Range($col + "3").Value = $row2[$col]
    $ws.Range($col + "4").Value = $row3[$col]
    $ws.Range($col + "2").Value = $row4[$col]
}
